# Applies the 2020-10-08 COVID19 time-series update: appends a repeated
# header row (757) plus 34 state/UT data rows (758:792) to the single
# worksheet, mirroring the layout used for every earlier daily block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 757: repeated header row. Copy the formatting (bold / boxed /
# centered) from the previous repeated header block at row 721 so the new
# block matches the existing pattern used throughout the sheet. ---
$ws.Range("A721:H721").Copy()
$ws.Range("A757:H757").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Row 757 header labels
$ws.Range("A757").Value = "States/UT"
$ws.Range("B757").Value = "Active Cases"
$ws.Range("C757").Value = "Active Cases Since Yesterday"
$ws.Range("D757").Value = "Recovered Cases"
$ws.Range("E757").Value = "Recovered Cases Since Yesterday"
$ws.Range("F757").Value = "Deceased Cases"
$ws.Range("G757").Value = "Deceased Cases Since Yesterday"
$ws.Range("H757").Value = "Date"

# Row 758
$ws.Range("A758").Value = "Andaman and Nicobar Islands"
$ws.Range("B758").Value = 185
$ws.Range("C758").Value = 5
$ws.Range("D758").Value = 3696
$ws.Range("E758").Value = 18
$ws.Range("F758").Value = 54
$ws.Range("G758").Value = 0
$ws.Range("H758").NumberFormat = "@"
$ws.Range("H758").Value = "08-10-2020"
$ws.Range("H758").Style = "Normal"

# Row 759
$ws.Range("A759").Value = "Andhra Pradesh"
$ws.Range("B759").Value = 49513
$ws.Range("C759").Value = -1263
$ws.Range("D759").Value = 678828
$ws.Range("E759").Value = 6349
$ws.Range("F759").Value = 6086
$ws.Range("G759").Value = 34
$ws.Range("H759").NumberFormat = "@"
$ws.Range("H759").Value = "08-10-2020"
$ws.Range("H759").Style = "Normal"

# Row 760
$ws.Range("A760").Value = "Arunachal Pradesh"
$ws.Range("B760").Value = 2850
$ws.Range("C760").Value = -172
$ws.Range("D760").Value = 8396
$ws.Range("E760").Value = 431
$ws.Range("F760").Value = 21
$ws.Range("G760").Value = 1
$ws.Range("H760").NumberFormat = "@"
$ws.Range("H760").Value = "08-10-2020"
$ws.Range("H760").Style = "Normal"

# Row 761
$ws.Range("A761").Value = "Assam"
$ws.Range("B761").Value = 31786
$ws.Range("C761").Value = -1261
$ws.Range("D761").Value = 157638
$ws.Range("E761").Value = 2561
$ws.Range("F761").Value = 785
$ws.Range("G761").Value = 7
$ws.Range("H761").NumberFormat = "@"
$ws.Range("H761").Value = "08-10-2020"
$ws.Range("H761").Style = "Normal"

# Row 762
$ws.Range("A762").Value = "Bihar"
$ws.Range("B762").Value = 11326
$ws.Range("C762").Value = -94
$ws.Range("D762").Value = 179732
$ws.Range("E762").Value = 1337
$ws.Range("F762").Value = 927
$ws.Range("G762").Value = 2
$ws.Range("H762").NumberFormat = "@"
$ws.Range("H762").Value = "08-10-2020"
$ws.Range("H762").Style = "Normal"

# Row 763
$ws.Range("A763").Value = "Chandigarh"
$ws.Range("B763").Value = 1448
$ws.Range("C763").Value = -44
$ws.Range("D763").Value = 11190
$ws.Range("E763").Value = 155
$ws.Range("F763").Value = 182
$ws.Range("G763").Value = 2
$ws.Range("H763").NumberFormat = "@"
$ws.Range("H763").Value = "08-10-2020"
$ws.Range("H763").Style = "Normal"

# Row 764
$ws.Range("A764").Value = "Chhattisgarh"
$ws.Range("B764").Value = 26777
$ws.Range("C764").Value = -461
$ws.Range("D764").Value = 103828
$ws.Range("E764").Value = 3277
$ws.Range("F764").Value = 1134
$ws.Range("G764").Value = 30
$ws.Range("H764").NumberFormat = "@"
$ws.Range("H764").Value = "08-10-2020"
$ws.Range("H764").Style = "Normal"

# Row 765
$ws.Range("A765").Value = "Dadra and Nagar Haveli and Daman and Diu"
$ws.Range("B765").Value = 108
$ws.Range("C765").Value = 7
$ws.Range("D765").Value = 3010
$ws.Range("E765").Value = 10
$ws.Range("F765").Value = 2
$ws.Range("G765").Value = 0
$ws.Range("H765").NumberFormat = "@"
$ws.Range("H765").Value = "08-10-2020"
$ws.Range("H765").Style = "Normal"

# Row 766
$ws.Range("A766").Value = "Delhi"
$ws.Range("B766").Value = 22186
$ws.Range("C766").Value = -534
$ws.Range("D766").Value = 270305
$ws.Range("E766").Value = 3370
$ws.Range("F766").Value = 5616
$ws.Range("G766").Value = 35
$ws.Range("H766").NumberFormat = "@"
$ws.Range("H766").Value = "08-10-2020"
$ws.Range("H766").Style = "Normal"

# Row 767
$ws.Range("A767").Value = "Goa"
$ws.Range("B767").Value = 4749
$ws.Range("C767").Value = 29
$ws.Range("D767").Value = 31444
$ws.Range("E767").Value = 394
$ws.Range("F767").Value = 477
$ws.Range("G767").Value = 9
$ws.Range("H767").NumberFormat = "@"
$ws.Range("H767").Value = "08-10-2020"
$ws.Range("H767").Style = "Normal"

# Row 768
$ws.Range("A768").Value = "Gujarat"
$ws.Range("B768").Value = 16485
$ws.Range("C768").Value = -85
$ws.Range("D768").Value = 126657
$ws.Range("E768").Value = 1546
$ws.Range("F768").Value = 3531
$ws.Range("G768").Value = 12
$ws.Range("H768").NumberFormat = "@"
$ws.Range("H768").Value = "08-10-2020"
$ws.Range("H768").Style = "Normal"

# Row 769
$ws.Range("A769").Value = "Haryana"
$ws.Range("B769").Value = 11029
$ws.Range("C769").Value = -291
$ws.Range("D769").Value = 124841
$ws.Range("E769").Value = 1555
$ws.Range("F769").Value = 1528
$ws.Range("G769").Value = 19
$ws.Range("H769").NumberFormat = "@"
$ws.Range("H769").Value = "08-10-2020"
$ws.Range("H769").Style = "Normal"

# Row 770
$ws.Range("A770").Value = "Himachal Pradesh"
$ws.Range("B770").Value = 2996
$ws.Range("C770").Value = -140
$ws.Range("D770").Value = 13338
$ws.Range("E770").Value = 420
$ws.Range("F770").Value = 231
$ws.Range("G770").Value = 2
$ws.Range("H770").NumberFormat = "@"
$ws.Range("H770").Value = "08-10-2020"
$ws.Range("H770").Style = "Normal"

# Row 771
$ws.Range("A771").Value = "Jammu and Kashmir"
$ws.Range("B771").Value = 12131
$ws.Range("C771").Value = -1581
$ws.Range("D771").Value = 67684
$ws.Range("E771").Value = 2188
$ws.Range("F771").Value = 1282
$ws.Range("G771").Value = 14
$ws.Range("H771").NumberFormat = "@"
$ws.Range("H771").Value = "08-10-2020"
$ws.Range("H771").Style = "Normal"

# Row 772
$ws.Range("A772").Value = "Jharkhand"
$ws.Range("B772").Value = 9759
$ws.Range("C772").Value = -268
$ws.Range("D772").Value = 79176
$ws.Range("E772").Value = 1087
$ws.Range("F772").Value = 767
$ws.Range("G772").Value = 10
$ws.Range("H772").NumberFormat = "@"
$ws.Range("H772").Value = "08-10-2020"
$ws.Range("H772").Style = "Normal"

# Row 773
$ws.Range("A773").Value = "Karnataka"
$ws.Range("B773").Value = 116172
$ws.Range("C773").Value = 1002
$ws.Range("D773").Value = 542906
$ws.Range("E773").Value = 9832
$ws.Range("F773").Value = 9574
$ws.Range("G773").Value = 113
$ws.Range("H773").NumberFormat = "@"
$ws.Range("H773").Value = "08-10-2020"
$ws.Range("H773").Style = "Normal"

# Row 774
$ws.Range("A774").Value = "Kerala"
$ws.Range("B774").Value = 92246
$ws.Range("C774").Value = 4423
$ws.Range("D774").Value = 160253
$ws.Range("E774").Value = 6161
$ws.Range("F774").Value = 906
$ws.Range("G774").Value = 22
$ws.Range("H774").NumberFormat = "@"
$ws.Range("H774").Value = "08-10-2020"
$ws.Range("H774").Style = "Normal"

# Row 775
$ws.Range("A775").Value = "Ladakh"
$ws.Range("B775").Value = 1228
$ws.Range("C775").Value = 33
$ws.Range("D775").Value = 3511
$ws.Range("E775").Value = 47
$ws.Range("F775").Value = 63
$ws.Range("G775").Value = 2
$ws.Range("H775").NumberFormat = "@"
$ws.Range("H775").Value = "08-10-2020"
$ws.Range("H775").Style = "Normal"

# Row 776
$ws.Range("A776").Value = "Madhya Pradesh"
$ws.Range("B776").Value = 17522
$ws.Range("C776").Value = -619
$ws.Range("D776").Value = 120267
$ws.Range("E776").Value = 2228
$ws.Range("F776").Value = 2518
$ws.Range("G776").Value = 30
$ws.Range("H776").NumberFormat = "@"
$ws.Range("H776").Value = "08-10-2020"
$ws.Range("H776").Style = "Normal"

# Row 777
$ws.Range("A777").Value = "Maharashtra"
$ws.Range("B777").Value = 244976
$ws.Range("C777").Value = -2492
$ws.Range("D777").Value = 1196441
$ws.Range("E777").Value = 16715
$ws.Range("F777").Value = 39072
$ws.Range("G777").Value = 355
$ws.Range("H777").NumberFormat = "@"
$ws.Range("H777").Value = "08-10-2020"
$ws.Range("H777").Style = "Normal"

# Row 778
$ws.Range("A778").Value = "Manipur"
$ws.Range("B778").Value = 2805
$ws.Range("C778").Value = 125
$ws.Range("D778").Value = 9604
$ws.Range("E778").Value = 122
$ws.Range("F778").Value = 80
$ws.Range("G778").Value = 2
$ws.Range("H778").NumberFormat = "@"
$ws.Range("H778").Value = "08-10-2020"
$ws.Range("H778").Style = "Normal"

# Row 779
$ws.Range("A779").Value = "Meghalaya"
$ws.Range("B779").Value = 2411
$ws.Range("C779").Value = 40
$ws.Range("D779").Value = 4694
$ws.Range("E779").Value = 88
$ws.Range("F779").Value = 60
$ws.Range("G779").Value = 0
$ws.Range("H779").NumberFormat = "@"
$ws.Range("H779").Value = "08-10-2020"
$ws.Range("H779").Style = "Normal"

# Row 780
$ws.Range("A780").Value = "Mizoram"
$ws.Range("B780").Value = 231
$ws.Range("C780").Value = -30
$ws.Range("D780").Value = 1919
$ws.Range("E780").Value = 32
$ws.Range("F780").Value = 0
$ws.Range("G780").Value = 0
$ws.Range("H780").NumberFormat = "@"
$ws.Range("H780").Value = "08-10-2020"
$ws.Range("H780").Style = "Normal"

# Row 781
$ws.Range("A781").Value = "Nagaland"
$ws.Range("B781").Value = 1200
$ws.Range("C781").Value = 15
$ws.Range("D781").Value = 5498
$ws.Range("E781").Value = 38
$ws.Range("F781").Value = 17
$ws.Range("G781").Value = 0
$ws.Range("H781").NumberFormat = "@"
$ws.Range("H781").Value = "08-10-2020"
$ws.Range("H781").Style = "Normal"

# Row 782
$ws.Range("A782").Value = "Odisha"
$ws.Range("B782").Value = 26368
$ws.Range("C782").Value = -478
$ws.Range("D782").Value = 213672
$ws.Range("E782").Value = 3455
$ws.Range("F782").Value = 958
$ws.Range("G782").Value = 18
$ws.Range("H782").NumberFormat = "@"
$ws.Range("H782").Value = "08-10-2020"
$ws.Range("H782").Style = "Normal"

# Row 783
$ws.Range("A783").Value = "Puducherry"
$ws.Range("B783").Value = 4680
$ws.Range("C783").Value = 158
$ws.Range("D783").Value = 24930
$ws.Range("E783").Value = 316
$ws.Range("F783").Value = 551
$ws.Range("G783").Value = 5
$ws.Range("H783").NumberFormat = "@"
$ws.Range("H783").Value = "08-10-2020"
$ws.Range("H783").Style = "Normal"

# Row 784
$ws.Range("A784").Value = "Punjab"
$ws.Range("B784").Value = 11563
$ws.Range("C784").Value = -419
$ws.Range("D784").Value = 105585
$ws.Range("E784").Value = 1230
$ws.Range("F784").Value = 3712
$ws.Range("G784").Value = 33
$ws.Range("H784").NumberFormat = "@"
$ws.Range("H784").Value = "08-10-2020"
$ws.Range("H784").Style = "Normal"

# Row 785
$ws.Range("A785").Value = "Rajasthan"
$ws.Range("B785").Value = 21351
$ws.Range("C785").Value = 57
$ws.Range("D785").Value = 127526
$ws.Range("E785").Value = 2078
$ws.Range("F785").Value = 1590
$ws.Range("G785").Value = 16
$ws.Range("H785").NumberFormat = "@"
$ws.Range("H785").Value = "08-10-2020"
$ws.Range("H785").Style = "Normal"

# Row 786
$ws.Range("A786").Value = "Sikkim"
$ws.Range("B786").Value = 570
$ws.Range("C786").Value = -10
$ws.Range("D786").Value = 2615
$ws.Range("E786").Value = 28
$ws.Range("F786").Value = 49
$ws.Range("G786").Value = 0
$ws.Range("H786").NumberFormat = "@"
$ws.Range("H786").Value = "08-10-2020"
$ws.Range("H786").Style = "Normal"

# Row 787
$ws.Range("A787").Value = "Tamil Nadu"
$ws.Range("B787").Value = 45135
$ws.Range("C787").Value = -144
$ws.Range("D787").Value = 580736
$ws.Range("E787").Value = 5524
$ws.Range("F787").Value = 9984
$ws.Range("G787").Value = 67
$ws.Range("H787").NumberFormat = "@"
$ws.Range("H787").Value = "08-10-2020"
$ws.Range("H787").Style = "Normal"

# Row 788
$ws.Range("A788").Value = "Telengana"
$ws.Range("B788").Value = 26368
$ws.Range("C788").Value = -183
$ws.Range("D788").Value = 179075
$ws.Range("E788").Value = 2067
$ws.Range("F788").Value = 1201
$ws.Range("G788").Value = 12
$ws.Range("H788").NumberFormat = "@"
$ws.Range("H788").Value = "08-10-2020"
$ws.Range("H788").Style = "Normal"

# Row 789
$ws.Range("A789").Value = "Tripura"
$ws.Range("B789").Value = 4389
$ws.Range("C789").Value = -232
$ws.Range("D789").Value = 23066
$ws.Range("E789").Value = 443
$ws.Range("F789").Value = 304
$ws.Range("G789").Value = 3
$ws.Range("H789").NumberFormat = "@"
$ws.Range("H789").Value = "08-10-2020"
$ws.Range("H789").Style = "Normal"

# Row 790
$ws.Range("A790").Value = "Uttarakhand"
$ws.Range("B790").Value = 8367
$ws.Range("C790").Value = -47
$ws.Range("D790").Value = 43904
$ws.Range("E790").Value = 666
$ws.Range("F790").Value = 688
$ws.Range("G790").Value = 11
$ws.Range("H790").NumberFormat = "@"
$ws.Range("H790").Value = "08-10-2020"
$ws.Range("H790").Style = "Normal"

# Row 791
$ws.Range("A791").Value = "Uttar Pradesh"
$ws.Range("B791").Value = 43154
$ws.Range("C791").Value = -877
$ws.Range("D791").Value = 374972
$ws.Range("E791").Value = 4219
$ws.Range("F791").Value = 6200
$ws.Range("G791").Value = 47
$ws.Range("H791").NumberFormat = "@"
$ws.Range("H791").Value = "08-10-2020"
$ws.Range("H791").Style = "Normal"

# Row 792
$ws.Range("A792").Value = "West Bengal"
$ws.Range("B792").Value = 28361
$ws.Range("C792").Value = 373
$ws.Range("D792").Value = 246767
$ws.Range("E792").Value = 3024
$ws.Range("F792").Value = 5376
$ws.Range("G792").Value = 58
$ws.Range("H792").NumberFormat = "@"
$ws.Range("H792").Value = "08-10-2020"
$ws.Range("H792").Style = "Normal"

